$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "FAPs"
$ws.Cells.Item(2,2).Value = "Efna2"
$ws.Cells.Item(2,3).Value = "Epha4"
$ws.Cells.Item(2,4).Value = "FAPs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 1.192082333333333
$ws.Cells.Item(2,8).Value = 3.576247
$ws.Cells.Item(2,9).Value = 0.2797939869571494
$ws.Cells.Item(2,10).Value = 0.2797939869571493
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 6.708176333333333
$ws.Cells.Item(2,14).Value = 20.124529
$ws.Cells.Item(2,15).Value = 0.4356329228871633
$ws.Cells.Item(2,16).Value = 0.4356329228871633
$ws.Cells.Item(2,17).Value = 7.996698495851444
$ws.Cells.Item(2,18).Value = 71.970286462663
$ws.Cells.Item(2,19).Value = 0.1218874723443958
$ws.Cells.Item(2,20).Value = 0.1218874723443958

# Row 3
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Efna2"
$ws.Cells.Item(3,3).Value = "Epha4"
$ws.Cells.Item(3,4).Value = "sCs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 1.192082333333333
$ws.Cells.Item(3,8).Value = 3.576247
$ws.Cells.Item(3,9).Value = 0.2797939869571494
$ws.Cells.Item(3,10).Value = 0.2797939869571493
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 6.789877333333333
$ws.Cells.Item(3,14).Value = 20.369632
$ws.Cells.Item(3,15).Value = 0.4409386339573907
$ws.Cells.Item(3,16).Value = 0.4409386339573907
$ws.Cells.Item(3,17).Value = 8.09409281456711
$ws.Cells.Item(3,18).Value = 72.846835331104
$ws.Cells.Item(3,19).Value = 0.1233719783983774
$ws.Cells.Item(3,20).Value = 0.1233719783983774

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Efna2"
$ws.Cells.Item(4,3).Value = "Epha4"
$ws.Cells.Item(4,4).Value = "ECs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 1.192082333333333
$ws.Cells.Item(4,8).Value = 3.576247
$ws.Cells.Item(4,9).Value = 0.2797939869571494
$ws.Cells.Item(4,10).Value = 0.2797939869571493
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 1.900636333333334
$ws.Cells.Item(4,14).Value = 5.701909000000001
$ws.Cells.Item(4,15).Value = 0.1234284431554459
$ws.Cells.Item(4,16).Value = 0.1234284431554459
$ws.Cells.Item(4,17).Value = 2.265714995058111
$ws.Cells.Item(4,18).Value = 20.391434955523
$ws.Cells.Item(4,19).Value = 0.03453453621437608
$ws.Cells.Item(4,20).Value = 0.03453453621437608

# Row 5
$ws.Cells.Item(5,1).Value = "sCs"
$ws.Cells.Item(5,2).Value = "Efna2"
$ws.Cells.Item(5,3).Value = "Epha4"
$ws.Cells.Item(5,4).Value = "FAPs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 1.843761666666667
$ws.Cells.Item(5,8).Value = 5.531285
$ws.Cells.Item(5,9).Value = 0.4327498305196134
$ws.Cells.Item(5,10).Value = 0.4327498305196134
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 6.708176333333333
$ws.Cells.Item(5,14).Value = 20.124529
$ws.Cells.Item(5,15).Value = 0.4356329228871633
$ws.Cells.Item(5,16).Value = 0.4356329228871633
$ws.Cells.Item(5,17).Value = 12.36827837664056
$ws.Cells.Item(5,18).Value = 111.314505389765
$ws.Cells.Item(5,19).Value = 0.1885200735481838
$ws.Cells.Item(5,20).Value = 0.1885200735481837

# Row 6
$ws.Cells.Item(6,1).Value = "sCs"
$ws.Cells.Item(6,2).Value = "Efna2"
$ws.Cells.Item(6,3).Value = "Epha4"
$ws.Cells.Item(6,4).Value = "sCs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 1.843761666666667
$ws.Cells.Item(6,8).Value = 5.531285
$ws.Cells.Item(6,9).Value = 0.4327498305196134
$ws.Cells.Item(6,10).Value = 0.4327498305196134
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 6.789877333333333
$ws.Cells.Item(6,14).Value = 20.369632
$ws.Cells.Item(6,15).Value = 0.4409386339573907
$ws.Cells.Item(6,16).Value = 0.4409386339573907
$ws.Cells.Item(6,17).Value = 12.51891554856889
$ws.Cells.Item(6,18).Value = 112.67023993712
$ws.Cells.Item(6,19).Value = 0.1908161191146107
$ws.Cells.Item(6,20).Value = 0.1908161191146107

# Row 7
$ws.Cells.Item(7,1).Value = "sCs"
$ws.Cells.Item(7,2).Value = "Efna2"
$ws.Cells.Item(7,3).Value = "Epha4"
$ws.Cells.Item(7,4).Value = "ECs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 1.843761666666667
$ws.Cells.Item(7,8).Value = 5.531285
$ws.Cells.Item(7,9).Value = 0.4327498305196134
$ws.Cells.Item(7,10).Value = 0.4327498305196134
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.900636333333334
$ws.Cells.Item(7,14).Value = 5.701909000000001
$ws.Cells.Item(7,15).Value = 0.1234284431554459
$ws.Cells.Item(7,16).Value = 0.1234284431554459
$ws.Cells.Item(7,17).Value = 3.50432041367389
$ws.Cells.Item(7,18).Value = 31.53888372306501
$ws.Cells.Item(7,19).Value = 0.05341363785681896
$ws.Cells.Item(7,20).Value = 0.05341363785681895

# Row 8
$ws.Cells.Item(8,1).Value = "ECs"
$ws.Cells.Item(8,2).Value = "Efna2"
$ws.Cells.Item(8,3).Value = "Epha4"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 1.224727666666667
$ws.Cells.Item(8,8).Value = 3.674183
$ws.Cells.Item(8,9).Value = 0.2874561825232373
$ws.Cells.Item(8,10).Value = 0.2874561825232373
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 6.708176333333333
$ws.Cells.Item(8,14).Value = 20.124529
$ws.Cells.Item(8,15).Value = 0.4356329228871633
$ws.Cells.Item(8,16).Value = 0.4356329228871633
$ws.Cells.Item(8,17).Value = 8.215689148311888
$ws.Cells.Item(8,18).Value = 73.941202334807
$ws.Cells.Item(8,19).Value = 0.1252253769945838
$ws.Cells.Item(8,20).Value = 0.1252253769945838

# Row 9
$ws.Cells.Item(9,1).Value = "ECs"
$ws.Cells.Item(9,2).Value = "Efna2"
$ws.Cells.Item(9,3).Value = "Epha4"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 1.224727666666667
$ws.Cells.Item(9,8).Value = 3.674183
$ws.Cells.Item(9,9).Value = 0.2874561825232373
$ws.Cells.Item(9,10).Value = 0.2874561825232373
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 6.789877333333333
$ws.Cells.Item(9,14).Value = 20.369632
$ws.Cells.Item(9,15).Value = 0.4409386339573907
$ws.Cells.Item(9,16).Value = 0.4409386339573907
$ws.Cells.Item(9,17).Value = 8.315750623406222
$ws.Cells.Item(9,18).Value = 74.84175561065601
$ws.Cells.Item(9,19).Value = 0.1267505364444026
$ws.Cells.Item(9,20).Value = 0.1267505364444026

# Row 10
$ws.Cells.Item(10,1).Value = "ECs"
$ws.Cells.Item(10,2).Value = "Efna2"
$ws.Cells.Item(10,3).Value = "Epha4"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 1.224727666666667
$ws.Cells.Item(10,8).Value = 3.674183
$ws.Cells.Item(10,9).Value = 0.2874561825232373
$ws.Cells.Item(10,10).Value = 0.2874561825232373
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 1.900636333333334
$ws.Cells.Item(10,14).Value = 5.701909000000001
$ws.Cells.Item(10,15).Value = 0.1234284431554459
$ws.Cells.Item(10,16).Value = 0.1234284431554459
$ws.Cells.Item(10,17).Value = 2.327761901705223
$ws.Cells.Item(10,18).Value = 20.949857115347
$ws.Cells.Item(10,19).Value = 0.03548026908425088
$ws.Cells.Item(10,20).Value = 0.03548026908425088
